$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add the new data points in rows 3 and 4 (mirrors row 2's G/L cells) ---
# Copy number-format / font from the existing row-2 cells so the new cells
# pick up matching formatting, then write the values.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("G4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("L2").Copy() | Out-Null
$ws.Range("L3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("L4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

$ws.Range("G3").Value = 210
$ws.Range("G4").Value = 210
$ws.Range("L3").Value = 9
$ws.Range("L4").Value = 9

# --- Update the sheet's view / selection state ---
$ws.Activate()
$ws.Range("J4").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
